# Update recomputed TPM-derived statistics in the LR-pairs sheet (Tnc-Itgav)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 7.759382333333334
$ws.Range("H2").Value = 23.278147
$ws.Range("I2").Value = 0.03884312320086166
$ws.Range("J2").Value = 0.04014625174390325
$ws.Range("M2").Value = 13.89934866666667
$ws.Range("N2").Value = 41.69804600000001
$ws.Range("O2").Value = 0.04853507553134179
$ws.Range("P2").Value = 0.04999273878390351
$ws.Range("Q2").Value = 107.8503604889736
$ws.Range("R2").Value = 970.6532444007621
$ws.Range("S2").Value = 0.001885253918427035
$ws.Range("T2").Value = 0.002007021076585786
# Row 3
$ws.Range("G3").Value = 7.759382333333334
$ws.Range("H3").Value = 23.278147
$ws.Range("I3").Value = 0.03884312320086166
$ws.Range("J3").Value = 0.04014625174390325
$ws.Range("O3").Value = 0.245697991654417
$ws.Range("P3").Value = 0.253077086664408
$ws.Range("Q3").Value = 545.9683884541206
$ws.Range("R3").Value = 4913.715496087085
$ws.Range("S3").Value = 0.009543677360036801
$ws.Range("T3").Value = 0.01016009643184294
# Row 4
$ws.Range("G4").Value = 7.759382333333334
$ws.Range("H4").Value = 23.278147
$ws.Range("I4").Value = 0.03884312320086166
$ws.Range("J4").Value = 0.04014625174390325
$ws.Range("M4").Value = 82.007665
$ws.Range("N4").Value = 246.022995
$ws.Range("O4").Value = 0.2863622109480123
$ws.Range("P4").Value = 0.2949625822722868
$ws.Range("Q4").Value = 636.3288269989184
$ws.Range("R4").Value = 5726.959442990265
$ws.Range("S4").Value = 0.01112320263992478
$ws.Range("T4").Value = 0.011841642082935
# Row 5
$ws.Range("G5").Value = 7.759382333333334
$ws.Range("H5").Value = 23.278147
$ws.Range("I5").Value = 0.03884312320086166
$ws.Range("J5").Value = 0.04014625174390325
$ws.Range("M5").Value = 25.0501465
$ws.Range("N5").Value = 50.100293
$ws.Range("O5").Value = 0.0874724982879541
$ws.Range("P5").Value = 0.06006638442832619
$ws.Range("Q5").Value = 194.3736641995118
$ws.Range("R5").Value = 1166.241985197071
$ws.Range("S5").Value = 0.003397705027686162
$ws.Range("T5").Value = 0.002411440190605653
# Row 6
$ws.Range("G6").Value = 7.759382333333334
$ws.Range("H6").Value = 23.278147
$ws.Range("I6").Value = 0.03884312320086166
$ws.Range("J6").Value = 0.04014625174390325
$ws.Range("M6").Value = 95.05788666666668
$ws.Range("N6").Value = 285.17366
$ws.Range("O6").Value = 0.3319322235782747
$ws.Range("P6").Value = 0.3419012078510756
$ws.Range("Q6").Value = 737.5904864453357
$ws.Range("R6").Value = 6638.314378008021
$ws.Range("S6").Value = 0.01289328425478688
$ws.Range("T6").Value = 0.01372605196193387
# Row 7
$ws.Range("I7").Value = 0.8631909770948131
$ws.Range("J7").Value = 0.8921497401307179
$ws.Range("M7").Value = 13.89934866666667
$ws.Range("N7").Value = 41.69804600000001
$ws.Range("O7").Value = 0.04853507553134179
$ws.Range("P7").Value = 0.04999273878390351
$ws.Range("Q7").Value = 2396.703724597506
$ws.Range("R7").Value = 21570.33352137755
$ws.Range("S7").Value = 0.04189503927126947
$ws.Range("T7").Value = 0.04460100891448238
# Row 8
$ws.Range("I8").Value = 0.8631909770948131
$ws.Range("J8").Value = 0.8921497401307179
$ws.Range("O8").Value = 0.245697991654417
$ws.Range("P8").Value = 0.253077086664408
$ws.Range("S8").Value = 0.2120842894864094
$ws.Range("T8").Value = 0.2257826571006908
# Row 9
$ws.Range("I9").Value = 0.8631909770948131
$ws.Range("J9").Value = 0.8921497401307179
$ws.Range("M9").Value = 82.007665
$ws.Range("N9").Value = 246.022995
$ws.Range("O9").Value = 0.2863622109480123
$ws.Range("P9").Value = 0.2949625822722868
$ws.Range("Q9").Value = 14140.81198080921
$ws.Range("R9").Value = 127267.3078272829
$ws.Range("S9").Value = 0.2471852766712458
$ws.Range("T9").Value = 0.2631507911225062
# Row 10
$ws.Range("I10").Value = 0.8631909770948131
$ws.Range("J10").Value = 0.8921497401307179
$ws.Range("M10").Value = 25.0501465
$ws.Range("N10").Value = 50.100293
$ws.Range("O10").Value = 0.0874724982879541
$ws.Range("P10").Value = 0.06006638442832619
$ws.Range("Q10").Value = 4319.467110156421
$ws.Range("R10").Value = 25916.80266093852
$ws.Range("S10").Value = 0.07550547126610346
$ws.Range("T10").Value = 0.05358820925832301
# Row 11
$ws.Range("I11").Value = 0.8631909770948131
$ws.Range("J11").Value = 0.8921497401307179
$ws.Range("M11").Value = 95.05788666666668
$ws.Range("N11").Value = 285.17366
$ws.Range("O11").Value = 0.3319322235782747
$ws.Range("P11").Value = 0.3419012078510756
$ws.Range("Q11").Value = 16391.09835216506
$ws.Range("R11").Value = 147519.8851694855
$ws.Range("S11").Value = 0.2865209003997849
$ws.Range("T11").Value = 0.3050270737347157
# Row 12
$ws.Range("G12").Value = 0.05240566666666666
$ws.Range("H12").Value = 0.157217
$ws.Range("I12").Value = 0.0002623404388789996
$ws.Range("J12").Value = 0.0002711415672571033
$ws.Range("M12").Value = 13.89934866666667
$ws.Range("N12").Value = 41.69804600000001
$ws.Range("O12").Value = 0.04853507553134179
$ws.Range("P12").Value = 0.04999273878390351
$ws.Range("Q12").Value = 0.7284046331091112
$ws.Range("R12").Value = 6.555641697982001
$ws.Range("S12").Value = 0.0000127327130159176
$ws.Range("T12").Value = 0.00001355510954534257
# Row 13
$ws.Range("G13").Value = 0.05240566666666666
$ws.Range("H13").Value = 0.157217
$ws.Range("I13").Value = 0.0002623404388789996
$ws.Range("J13").Value = 0.0002711415672571033
$ws.Range("O13").Value = 0.245697991654417
$ws.Range("P13").Value = 0.253077086664408
$ws.Range("Q13").Value = 3.687385947326111
$ws.Range("R13").Value = 33.186473525935
$ws.Range("S13").Value = 0.00006445651896230853
$ws.Range("T13").Value = 0.00006861971791504934
# Row 14
$ws.Range("G14").Value = 0.05240566666666666
$ws.Range("H14").Value = 0.157217
$ws.Range("I14").Value = 0.0002623404388789996
$ws.Range("J14").Value = 0.0002711415672571033
$ws.Range("M14").Value = 82.007665
$ws.Range("N14").Value = 246.022995
$ws.Range("O14").Value = 0.2863622109480123
$ws.Range("P14").Value = 0.2949625822722868
$ws.Range("Q14").Value = 4.297666356101667
$ws.Range("R14").Value = 38.678997204915
$ws.Range("S14").Value = 0.00007512438809846222
$ws.Range("T14").Value = 0.00007997661683951012
# Row 15
$ws.Range("G15").Value = 0.05240566666666666
$ws.Range("H15").Value = 0.157217
$ws.Range("I15").Value = 0.0002623404388789996
$ws.Range("J15").Value = 0.0002711415672571033
$ws.Range("M15").Value = 25.0501465
$ws.Range("N15").Value = 50.100293
$ws.Range("O15").Value = 0.0874724982879541
$ws.Range("P15").Value = 0.06006638442832619
$ws.Range("Q15").Value = 1.312769627430167
$ws.Range("R15").Value = 7.876617764581
$ws.Range("S15").Value = 0.00002294757359070442
$ws.Range("T15").Value = 0.00001628649361336403
# Row 16
$ws.Range("G16").Value = 0.05240566666666666
$ws.Range("H16").Value = 0.157217
$ws.Range("I16").Value = 0.0002623404388789996
$ws.Range("J16").Value = 0.0002711415672571033
$ws.Range("M16").Value = 95.05788666666668
$ws.Range("N16").Value = 285.17366
$ws.Range("O16").Value = 0.3319322235782747
$ws.Range("P16").Value = 0.3419012078510756
$ws.Range("Q16").Value = 4.981571922691111
$ws.Range("R16").Value = 44.83414730422
$ws.Range("S16").Value = 0.00008707924521160678
$ws.Range("T16").Value = 0.00009270362934383728
# Row 17
$ws.Range("G17").Value = 19.452549
$ws.Range("H17").Value = 38.905098
$ws.Range("I17").Value = 0.09737859599105524
$ws.Range("J17").Value = 0.06709700125311635
$ws.Range("M17").Value = 13.89934866666667
$ws.Range("N17").Value = 41.69804600000001
$ws.Range("O17").Value = 0.04853507553134179
$ws.Range("P17").Value = 0.04999273878390351
$ws.Range("Q17").Value = 270.377761006418
$ws.Range("R17").Value = 1622.266566038508
$ws.Range("S17").Value = 0.004726277511561882
$ws.Range("T17").Value = 0.003354362856830292
# Row 18
$ws.Range("G18").Value = 19.452549
$ws.Range("H18").Value = 38.905098
$ws.Range("I18").Value = 0.09737859599105524
$ws.Range("J18").Value = 0.06709700125311635
$ws.Range("O18").Value = 0.245697991654417
$ws.Range("P18").Value = 0.253077086664408
$ws.Range("Q18").Value = 1368.727093551065
$ws.Range("R18").Value = 8212.362561306389
$ws.Range("S18").Value = 0.02392572546512913
$ws.Range("T18").Value = 0.01698071360105682
# Row 19
$ws.Range("G19").Value = 19.452549
$ws.Range("H19").Value = 38.905098
$ws.Range("I19").Value = 0.09737859599105524
$ws.Range("J19").Value = 0.06709700125311635
$ws.Range("M19").Value = 82.007665
$ws.Range("N19").Value = 246.022995
$ws.Range("O19").Value = 0.2863622109480123
$ws.Range("P19").Value = 0.2949625822722868
$ws.Range("Q19").Value = 1595.258121788085
$ws.Range("R19").Value = 9571.548730728509
$ws.Range("S19").Value = 0.02788555004701183
$ws.Range("T19").Value = 0.01979110475234606
# Row 20
$ws.Range("G20").Value = 19.452549
$ws.Range("H20").Value = 38.905098
$ws.Range("I20").Value = 0.09737859599105524
$ws.Range("J20").Value = 0.06709700125311635
$ws.Range("M20").Value = 25.0501465
$ws.Range("N20").Value = 50.100293
$ws.Range("O20").Value = 0.0874724982879541
$ws.Range("P20").Value = 0.06006638442832619
$ws.Range("Q20").Value = 487.2892022484285
$ws.Range("R20").Value = 1949.156808993714
$ws.Range("S20").Value = 0.008517949071110953
$ws.Range("T20").Value = 0.004030274271257571
# Row 21
$ws.Range("G21").Value = 19.452549
$ws.Range("H21").Value = 38.905098
$ws.Range("I21").Value = 0.09737859599105524
$ws.Range("J21").Value = 0.06709700125311635
$ws.Range("M21").Value = 95.05788666666668
$ws.Range("N21").Value = 285.17366
$ws.Range("O21").Value = 0.3319322235782747
$ws.Range("P21").Value = 0.3419012078510756
$ws.Range("Q21").Value = 1849.11819821978
$ws.Range("R21").Value = 11094.70918931868
$ws.Range("S21").Value = 0.03232309389624143
$ws.Range("T21").Value = 0.02294054577162562
# Row 22
$ws.Range("G22").Value = 0.06491533333333334
$ws.Range("H22").Value = 0.194746
$ws.Range("I22").Value = 0.0003249632743909987
$ws.Range("J22").Value = 0.0003358653050055137
$ws.Range("M22").Value = 13.89934866666667
$ws.Range("N22").Value = 41.69804600000001
$ws.Range("O22").Value = 0.04853507553134179
$ws.Range("P22").Value = 0.04999273878390351
$ws.Range("Q22").Value = 0.9022808518128891
$ws.Range("R22").Value = 8.120527666316001
$ws.Range("S22").Value = 0.00001577211706747927
$ws.Range("T22").Value = 0.00001679082645971672
# Row 23
$ws.Range("G23").Value = 0.06491533333333334
$ws.Range("H23").Value = 0.194746
$ws.Range("I23").Value = 0.0003249632743909987
$ws.Range("J23").Value = 0.0003358653050055137
$ws.Range("O23").Value = 0.245697991654417
$ws.Range("P23").Value = 0.253077086664408
$ws.Range("Q23").Value = 4.567595512558889
$ws.Range("R23").Value = 41.10835961303
$ws.Range("S23").Value = 0.00007984282387931164
$ws.Range("T23").Value = 0.0000849998129024482
# Row 24
$ws.Range("G24").Value = 0.06491533333333334
$ws.Range("H24").Value = 0.194746
$ws.Range("I24").Value = 0.0003249632743909987
$ws.Range("J24").Value = 0.0003358653050055137
$ws.Range("M24").Value = 82.007665
$ws.Range("N24").Value = 246.022995
$ws.Range("O24").Value = 0.2863622109480123
$ws.Range("P24").Value = 0.2949625822722868
$ws.Range("Q24").Value = 5.323554909363334
$ws.Range("R24").Value = 47.91199418427
$ws.Range("S24").Value = 0.000093057201731512
$ws.Range("T24").Value = 0.00009906769766009551
# Row 25
$ws.Range("G25").Value = 0.06491533333333334
$ws.Range("H25").Value = 0.194746
$ws.Range("I25").Value = 0.0003249632743909987
$ws.Range("J25").Value = 0.0003358653050055137
$ws.Range("M25").Value = 25.0501465
$ws.Range("N25").Value = 50.100293
$ws.Range("O25").Value = 0.0874724982879541
$ws.Range("P25").Value = 0.06006638442832619
$ws.Range("Q25").Value = 1.626138610096334
$ws.Range("R25").Value = 9.756831660577999
$ws.Range("S25").Value = 0.0000284253494628146
$ws.Range("T25").Value = 0.00002017421452659821
# Row 26
$ws.Range("G26").Value = 0.06491533333333334
$ws.Range("H26").Value = 0.194746
$ws.Range("I26").Value = 0.0003249632743909987
$ws.Range("J26").Value = 0.0003358653050055137
$ws.Range("M26").Value = 95.05788666666668
$ws.Range("N26").Value = 285.17366
$ws.Range("O26").Value = 0.3319322235782747
$ws.Range("P26").Value = 0.3419012078510756
$ws.Range("Q26").Value = 6.17071439892889
$ws.Range("R26").Value = 55.53642959036001
$ws.Range("S26").Value = 0.0001078657822498812
$ws.Range("T26").Value = 0.000114832753456655
